$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window tab ratio (995 -> 500) ---
$excel.Windows.Item(1).TabRatio = 500

# --- Currency text fix: EUR -> USD (shared string, updates all K/O cells) ---
$ws.Cells.Replace("EUR", "USD")

# --- Data fix on row 10: Gross Debit (GC) 666 -> 1598 ---
$ws.Range("M10").Value = 1598

# --- Date/time number format: uppercase tokens -> lowercase tokens ---
$ws.Range("G5:G25").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# --- Column T (20) gets an explicit custom width ---
$ws.Columns.Item(20).ColumnWidth = 16.6367

# --- Selection / view moves to L9 ---
$ws.Range("L9").Select()
